$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet 1 (DQ_Report)
# ---------------------------------------------------------------------------

# Header row -> new English column titles, bold
$ws1.Range("A1").Value = "Patient ID"
$ws1.Range("B1").Value = "Admission ID"
$ws1.Range("C1").Value = "ICD_Primary Code"
$ws1.Range("D1").Value = "Orphacode"
$ws1.Range("E1").Value = "DQ_Violations"
$ws1.Range("A1:E1").Font.Bold = $true

# Update the implausible-birthdate max-age message (105 -> 130)
$ws1.Range("E3").Value = "Implausible birthdate 1877-12-01 maximal age 130. "

# Move the "Following items are missing..." message up into row 11 (E11),
# then delete the now-empty row 12
$ws1.Range("E11").Value = "Following items are missing:  Kontakt_Klasse , Fall_Status , DiagnoseRolle"
$ws1.Rows.Item(12).Delete()

# Column widths
$ws1.Columns.Item(1).ColumnWidth = 26.71
$ws1.Columns.Item(2).ColumnWidth = 26.71
$ws1.Columns.Item(3).ColumnWidth = 16.71
$ws1.Columns.Item(4).ColumnWidth = 9.71
$ws1.Columns.Item(5).ColumnWidth = 101.71

# ---------------------------------------------------------------------------
# Sheet 2 (DQ_Metrics)
# ---------------------------------------------------------------------------

# Header row 1 -- columns A-K unchanged; L.. onward re-laid out with new
# metric columns inserted, existing ones shifted right.
$ws2.Range("A1").Value = "inst_id"
$ws2.Range("B1").Value = "report_year"
$ws2.Range("C1").Value = "item_completeness_rate"
$ws2.Range("D1").Value = "value_completeness_rate"
$ws2.Range("E1").Value = "orphaCoding_completeness_rate"
$ws2.Range("F1").Value = "range_plausibility_rate"
$ws2.Range("G1").Value = "orphaCoding_plausibility_rate"
$ws2.Range("H1").Value = "rdCase_unambiguity_rate"
$ws2.Range("I1").Value = "rdCase_dissimilarity_rate"
$ws2.Range("J1").Value = "case_no_py_ipat"
$ws2.Range("K1").Value = "case_no_py"
$ws2.Range("L1").Value = "missing_item_no_py"
$ws2.Range("M1").Value = "missing_value_no_py"
$ws2.Range("N1").Value = "outlier_no_py"
$ws2.Range("O1").Value = "orphaMissing_no_py"
$ws2.Range("P1").Value = "implausible_codeLink_no_py"
$ws2.Range("Q1").Value = "ambiguous_rdCase_no_py"
$ws2.Range("R1").Value = "duplicateRdCase_no_py"
$ws2.Range("S1").Value = "rdCase_no_py"
$ws2.Range("T1").Value = "mxCases_no_py"
$ws2.Range("U1").Value = "orphaCase_no_py"
$ws2.Range("V1").Value = "tracerCase_no_py"
$ws2.Range("W1").Value = "rdCase_rel_py_ipat"
$ws2.Range("X1").Value = "orphaCase_rel_py_ipat"
$ws2.Range("Y1").Value = "tracerCase_rel_py_ipat"
$ws2.Range("Z1").Value = "executionTime_inMin"
$ws2.Range("AA1").Value = "dateRef"
$ws2.Range("AB1").Value = "dataFormat"
$ws2.Range("AC1").Value = "diagnosesList"
$ws2.Range("AD1").Value = "encounterClass"
$ws2.Range("A1:AD1").Font.Bold = $true

# Row 2 data values
$ws2.Range("A2").Value = "260123451-Airolo"
$ws2.Range("B2").Value = 2020
$ws2.Range("C2").Value = 78.57
$ws2.Range("D2").Value = 99.03
$ws2.Range("E2").Value = 58.33
$ws2.Range("F2").Value = 99.46
$ws2.Range("G2").Value = 86.96
$ws2.Range("H2").Value = 88.89
$ws2.Range("I2").Value = 100
$ws2.Range("J2").Value = 997
$ws2.Range("K2").Value = 45
$ws2.Range("L2").Value = 3
$ws2.Range("M2").Value = 5
$ws2.Range("N2").Value = 1
$ws2.Range("O2").Value = 5
$ws2.Range("P2").Value = 3
$ws2.Range("Q2").Value = 3
$ws2.Range("R2").Value = 0
$ws2.Range("S2").Value = 27
$ws2.Range("T2").Value = 18
$ws2.Range("U2").Value = 23
$ws2.Range("V2").Value = 11
$ws2.Range("W2").Value = 2708
$ws2.Range("X2").Value = 2307
$ws2.Range("Y2").Value = 1103
$ws2.Range("Z2").Value = 0.04
$ws2.Range("AA2").Value = "Diagnosedatum"
$ws2.Range("AB2").Value = "FHIR"
$ws2.Range("AC2").Value = "v1"

# Column widths A..AD
$ws2.Columns.Item(1).ColumnWidth = 16.71
$ws2.Columns.Item(2).ColumnWidth = 11.71
$ws2.Columns.Item(3).ColumnWidth = 22.71
$ws2.Columns.Item(4).ColumnWidth = 23.71
$ws2.Columns.Item(5).ColumnWidth = 29.71
$ws2.Columns.Item(6).ColumnWidth = 23.71
$ws2.Columns.Item(7).ColumnWidth = 29.71
$ws2.Columns.Item(8).ColumnWidth = 23.71
$ws2.Columns.Item(9).ColumnWidth = 25.71
$ws2.Columns.Item(10).ColumnWidth = 15.71
$ws2.Columns.Item(11).ColumnWidth = 10.71
$ws2.Columns.Item(12).ColumnWidth = 18.71
$ws2.Columns.Item(13).ColumnWidth = 19.71
$ws2.Columns.Item(14).ColumnWidth = 13.71
$ws2.Columns.Item(15).ColumnWidth = 18.71
$ws2.Columns.Item(16).ColumnWidth = 26.71
$ws2.Columns.Item(17).ColumnWidth = 22.71
$ws2.Columns.Item(18).ColumnWidth = 21.71
$ws2.Columns.Item(19).ColumnWidth = 12.71
$ws2.Columns.Item(20).ColumnWidth = 13.71
$ws2.Columns.Item(21).ColumnWidth = 15.71
$ws2.Columns.Item(22).ColumnWidth = 16.71
$ws2.Columns.Item(23).ColumnWidth = 18.71
$ws2.Columns.Item(24).ColumnWidth = 21.71
$ws2.Columns.Item(25).ColumnWidth = 22.71
$ws2.Columns.Item(26).ColumnWidth = 19.71
$ws2.Columns.Item(27).ColumnWidth = 13.71
$ws2.Columns.Item(28).ColumnWidth = 10.71
$ws2.Columns.Item(29).ColumnWidth = 13.71
$ws2.Columns.Item(30).ColumnWidth = 14.71

Write-Host "edit complete"
